# data : case 1
# Update the 6 numeric values in A1:B3 (row 4 is untouched) and widen
# columns A and B slightly, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.038889402114658814
$ws.Range("B1").Value = -0.038889402144032012

$ws.Range("A2").Value = 0.0065005019915595061
$ws.Range("B2").Value = -0.0065005020338793369

$ws.Range("A3").Value = -0.062578075209440323
$ws.Range("B3").Value = 0.062578075175741293

# Target OOXML col widths are 14.7109375 (A) and 15.42578125 (B). This
# engine's COM ColumnWidth stores widths on a 1/6-character (whole-point)
# grid, so set the nearest representable ColumnWidth (xml_width - 5/6) for
# each: round(14.7109375*6)/6 - 5/6 and round(15.42578125*6)/6 - 5/6.
$ws.Columns.Item(1).ColumnWidth = 13.833333333333332
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666
